$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook logs a WhatsApp-style conversation; this change appends
# 16 newly-captured messages (rows 17-32) from the "Noah" thread, extending the
# sheet from A1:G16 to A1:G32. Row 16 (an existing "Noah" row) is used as a
# template so the new rows inherit the same cell formatting/typing (e.g. the
# Phone column stays text) before the per-row values are filled in.
$template = $ws.Range("A16:G16")

$rows = @(
    @(17, "2025-09-23 12:49:42", "Noah", 8450689526, "13052054965", "Hello", "media_files\photo_2025-09-23_16-49-42.jpg"),
    @(18, "2025-09-23 15:05:59", "Noah", 8450689526, "13052054965", "Hi", "media_files\photo_2025-09-23_19-05-59.jpg"),
    @(19, "2025-09-23 15:07:25", "Noah", 8450689526, "13052054965", "Hiii", "my-node-server/uploads/images\photo_2025-09-23_19-07-25.jpg"),
    @(20, "2025-09-23 15:20:01", "Noah", 8450689526, "13052054965", "Hello", "my-node-server/uploads/images\photo_2025-09-23_19-20-01.jpg"),
    @(21, "2025-09-23 15:30:16", "Noah", 8450689526, "13052054965", "Whats up", "my-node-server/uploads/images\photo_2025-09-23_19-30-16.jpg"),
    @(22, "2025-09-23 15:40:39", "Noah", 8450689526, "13052054965", "Yo", "my-node-server/public/uploads/images\photo_2025-09-23_19-40-39.jpg"),
    @(23, "2025-09-23 15:43:25", "Noah", 8450689526, "13052054965", "Hey", "my-node-server/public/uploads/images\photo_2025-09-23_19-43-25.jpg"),
    @(24, "2025-09-23 15:43:45", "Noah", 8450689526, "13052054965", "Yuh", "my-node-server/public/uploads/images\photo_2025-09-23_19-43-45.jpg"),
    @(25, "2025-09-23 15:58:02", "Noah", 8450689526, "13052054965", " ", "my-node-server/public/uploads/videos\video.mp4"),
    @(26, "2025-09-23 20:30:25", "Noah", 8450689526, "13052054965", "Yooooo", "my-node-server/public/uploads/images\photo_2025-09-24_00-30-25.jpg"),
    @(27, "2025-09-23 20:35:42", "Noah", 8450689526, "13052054965", "Hey man", "my-node-server/public/uploads/images\photo_2025-09-24_00-35-42.jpg"),
    @(28, "2025-09-23 20:43:43", "Noah", 8450689526, "13052054965", "Hey man", "my-node-server/public/uploads/images\photo_2025-09-24_00-43-43.jpg"),
    @(29, "2025-09-23 20:44:46", "Noah", 8450689526, "13052054965", "Hey man", "my-node-server/public/uploads/images\photo_2025-09-24_00-44-46.jpg"),
    @(30, "2025-09-23 20:48:08", "Noah", 8450689526, "13052054965", "What’s up man", "my-node-server/public/uploads/images\photo_2025-09-24_00-48-08.jpg"),
    @(31, "2025-09-23 20:54:17", "Noah", 8450689526, "13052054965", "Hey", "my-node-server/public/uploads/images\photo_2025-09-24_00-54-17.jpg"),
    @(32, "2025-09-23 20:57:01", "Noah", 8450689526, "13052054965", "What’s up man", "my-node-server/public/uploads/images\photo_2025-09-24_00-57-02.jpg")
)

foreach ($data in $rows) {
    $r = $data[0]

    # Copy the template row down first so the new row picks up the same
    # column types/styles (Sender Id stays numeric, Phone stays text, etc.)
    $template.Copy($ws.Range("A" + $r + ":G" + $r))

    $ws.Cells.Item($r, 1).Value = $data[1]
    $ws.Cells.Item($r, 2).Value = $data[2]
    $ws.Cells.Item($r, 3).Value = $data[3]
    # Column D (Phone) is left as copied from the template: every new row is
    # the same "13052054965" conversation, and re-assigning a digit-only
    # string here would make Excel reinterpret the cell as a number again.
    $ws.Cells.Item($r, 5).Value = $data[5]
    $ws.Cells.Item($r, 6).Value = $data[6]
}
